# logboek_jarno_touw.xlsx update
# - restyle a few existing "week4" cells (drop the ad-hoc duplicate styles 16/17
#   in favour of the already-existing equivalent styles 1/11)
# - add the new "finance" journal entries (rows 42,43,46,47,48,49)
# - move the active selection to G44

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-point the cells that were using the redundant styles (16 -> 1,
#    17 -> 11) back onto the canonical ones so the duplicate xf entries stop
#    being referenced. Clear() strips the current (duplicate) format; a
#    plain re-assignment of the text then lands back on the column's
#    canonical style (1) for G, while C39 needs the distinct yellow/bordered
#    style (11) restored explicitly via a format-only paste from a sibling
#    cell that already uses it.
# ---------------------------------------------------------------------------

$ws.Range("G37").Clear()
$ws.Range("G37").Value = "sales afdeling verder afgemaakt onder andere cutomers afdeling en projects afdeling"

$ws.Range("C39").Clear()
$ws.Range("C39").Value = "week4"
$ws.Range("C14").Copy()
$ws.Range("C39").PasteSpecial(-4122)

$ws.Range("G39").Clear()
$ws.Range("G39").Value = "sales afdeling verder afgemaakt onder andere cutomers afdeling en projects afdeling"

$ws.Range("G40").Clear()
$ws.Range("G40").Value = "header gefixed, admin panel frontend gemaakt"

# ---------------------------------------------------------------------------
# 2) New journal rows. Column E carries the date (style of row 40 / numFmt
#    14), column G carries the free-text note (column default style 1).
#    For the date cells we set the serial number first and only then paste
#    the number-format from an existing date cell, which reuses the existing
#    cellXfs entry instead of minting a new custom numFmt.
# ---------------------------------------------------------------------------

function Set-LogDate($cellAddr, $serial) {
    $ws.Range($cellAddr).Value = $serial
    $ws.Range("E40").Copy()
    $ws.Range($cellAddr).PasteSpecial(-4122)
}

Set-LogDate "E42" 43013
Set-LogDate "E43" 43014
Set-LogDate "E46" 43017
Set-LogDate "E47" 43018
Set-LogDate "E48" 43020
Set-LogDate "E49" 43021

# Text notes - the order below matters: it reproduces the exact order in
# which the new shared strings were first introduced.
$ws.Range("G42").Value = "admin panel verder gemaakt(front-end af)"
$ws.Range("G48").Value = "finance back end"
$ws.Range("G49").Value = "finance back end project info en klant en project inactief zetten"
$ws.Range("G46").Value = "begonnen aan finance back end"
$ws.Range("G47").Value = "begonnen aan finance back end"
$ws.Range("G43").Value = "begonnen aan finance front end/back end"

# ---------------------------------------------------------------------------
# 3) Leave the selection on G44, matching where the author ended up.
# ---------------------------------------------------------------------------
$ws.Range("G44").Select()
